# Update the bounding-box coordinates on the "Mapping" sheet.
# The shapefiles were reprojected to WGS 84 (EPSG 4269), so the reserve
# (Res_Bounding_Box) and sentinel-site (SK_Bounding_Box) longitude/latitude
# corners move slightly. Only columns A (Res_Bounding_Box) and B
# (SK_Bounding_Box) in rows 2-5 change; everything else on the sheet is
# unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -80.7155
$ws.Range("B2").Value = -80.6925

$ws.Range("A3").Value = 32.2987
$ws.Range("B3").Value = 32.318

$ws.Range("A4").Value = -80.1588
$ws.Range("B4").Value = -80.1818

$ws.Range("A5").Value = 32.7158
$ws.Range("B5").Value = 32.6965
